# Modify calculation method for saving percentage.
#
# "cheapest_price_saving" (column I) used to be computed as
#     cheapest_price / cheapest_price_baseline
# and "fastest_runtime_saving" (column P) used to be computed as
#     fastest_runtime / fastest_runtime_baseline
#
# Both are corrected to be an actual *saving* percentage, i.e. how much
# less than the baseline was spent/used:
#     (baseline - value) / baseline
#
# This touches every data row (2-33) of columns I and P on all six
# worksheets. The new text is written through a formula + copy /
# paste-special(values) round trip so the result lands back in the cell
# as a literal percentage string (matching how the original data was
# authored) instead of Excel auto-converting a "41.62%"-looking string
# into a percentage-formatted number.

$wb = $excel.ActiveWorkbook

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    for ($r = 2; $r -le 33; $r++) {
        $cheapestPrice = $ws.Cells.Item($r, 3).Value2()
        $cheapestBaseline = $ws.Cells.Item($r, 8).Value2()
        $fastestRuntime = $ws.Cells.Item($r, 11).Value2()
        $fastestBaseline = $ws.Cells.Item($r, 15).Value2()

        $cheapestSaving = ($cheapestBaseline - $cheapestPrice) / $cheapestBaseline * 100
        $fastestSaving = ($fastestBaseline - $fastestRuntime) / $fastestBaseline * 100

        $cheapestText = "{0:N2}%" -f $cheapestSaving
        $fastestText = "{0:N2}%" -f $fastestSaving

        $ws.Cells.Item($r, 9).Formula = '="' + $cheapestText + '"'
        $ws.Cells.Item($r, 16).Formula = '="' + $fastestText + '"'
    }

    $iRange = $ws.Range("I2:I33")
    $iRange.Copy()
    $iRange.PasteSpecial(-4163)

    $pRange = $ws.Range("P2:P33")
    $pRange.Copy()
    $pRange.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
